# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
# Rows 172-174 (match ids B172/B173/B174 etc.) were re-ordered: the data that
# used to live in row 173 now lives in row 172, row 174's old data now lives
# in row 173, and row 172's old data now lives in row 174. Columns A, C and D
# are identical across these three rows, so only B and E:AD need updating.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 172 (now holds what used to be row 173's match: FK Kapaz vs Zira IK)
$ws.Range("B172").Value = 7157321
$ws.Range("E172").Value = "FK Kapaz"
$ws.Range("F172").Value = "Zira IK"
$ws.Range("H172").Value = 1
$ws.Range("M172").Value = 3.4
$ws.Range("N172").Value = 1.666
$ws.Range("O172").Value = 5
$ws.Range("P172").Value = 3.6
$ws.Range("Q172").Value = 1.533
$ws.Range("R172").Value = 1
$ws.Range("S172").Value = 1.825
$ws.Range("T172").Value = 1.975
$ws.Range("U172").Value = 2.5
$ws.Range("V172").Value = 1.875
$ws.Range("W172").Value = 1.925
$ws.Range("Z172").Value = 0.5329999999999999
$ws.Range("AA172").Value = 0
$ws.Range("AB172").Value = 0
$ws.Range("AC172").Value = -1
$ws.Range("AD172").Value = 0.925

# Row 173 (now holds what used to be row 174's match: Sabah vs FK Sumqayit)
$ws.Range("B173").Value = 7153759
$ws.Range("E173").Value = "Sabah"
$ws.Range("F173").Value = "FK Sumqayit"
$ws.Range("G173").Value = 2
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 1
$ws.Range("K173").Value = "H"
$ws.Range("L173").Value = 1.833
$ws.Range("M173").Value = 3.3
$ws.Range("N173").Value = 3.8
$ws.Range("O173").Value = 2.15
$ws.Range("P173").Value = 3.2
$ws.Range("Q173").Value = 3.1
$ws.Range("R173").Value = -0.25
$ws.Range("S173").Value = 1.9
$ws.Range("T173").Value = 1.9
$ws.Range("U173").Value = 2.25
$ws.Range("V173").Value = 1.8
$ws.Range("W173").Value = 2
$ws.Range("X173").Value = 1.15
$ws.Range("Z173").Value = -1
$ws.Range("AA173").Value = 0.8999999999999999
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = -0.5
$ws.Range("AD173").Value = 0.5

# Row 174 (now holds what used to be row 172's match: Sabail FC vs Neftchi Baku)
$ws.Range("B174").Value = 7158118
$ws.Range("E174").Value = "Sabail FC"
$ws.Range("F174").Value = "Neftchi Baku"
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 3
$ws.Range("I174").Value = 0
$ws.Range("K174").Value = "A"
$ws.Range("L174").Value = 4.5
$ws.Range("M174").Value = 3.7
$ws.Range("N174").Value = 1.6
$ws.Range("O174").Value = 4.2
$ws.Range("P174").Value = 3.5
$ws.Range("Q174").Value = 1.7
$ws.Range("R174").Value = 0.75
$ws.Range("U174").Value = 2.75
$ws.Range("V174").Value = 1.825
$ws.Range("W174").Value = 1.975
$ws.Range("X174").Value = -1
$ws.Range("Z174").Value = 0.7
$ws.Range("AA174").Value = -1
$ws.Range("AB174").Value = 0.8999999999999999
$ws.Range("AC174").Value = 0.4125
$ws.Range("AD174").Value = -0.5
